# Weekly update: "Fruta / hortaliza, semanal"
# Two new rows of price data are inserted at the top of the
# "Vega Monumental Concepción - Nectarín" price list (currently rows 488-489),
# pushing the existing rows 488-514 down to 490-516.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 488 - this shifts
# old rows 488:514 down to 490:516, matching the diff's renumbering.
$ws.Rows("488:489").Insert()

# --- New row 488: Artic Star, Primera ---
$ws.Cells.Item(488, 1).Value = 11
$ws.Cells.Item(488, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(488, 3).Value = "Bíobío"
$ws.Cells.Item(488, 4).Value = 45267
$ws.Cells.Item(488, 5).Value = 8
$ws.Cells.Item(488, 6).Value = "Fruta"
$ws.Cells.Item(488, 7).Value = 100103
$ws.Cells.Item(488, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(488, 9).Value = 100103006
$ws.Cells.Item(488, 10).Value = "Nectarín"
$ws.Cells.Item(488, 11).Value = "Artic Star"
$ws.Cells.Item(488, 12).Value = "Primera"
$ws.Cells.Item(488, 13).Value = 300
$ws.Cells.Item(488, 14).Value = 17000
$ws.Cells.Item(488, 15).Value = 17000
$ws.Cells.Item(488, 16).Value = 17000
$ws.Cells.Item(488, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(488, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(488, 19).Value = 944
$ws.Cells.Item(488, 20).Value = 18

# --- New row 489: Early Glo, Primera ---
$ws.Cells.Item(489, 1).Value = 11
$ws.Cells.Item(489, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(489, 3).Value = "Bíobío"
$ws.Cells.Item(489, 4).Value = 45267
$ws.Cells.Item(489, 5).Value = 8
$ws.Cells.Item(489, 6).Value = "Fruta"
$ws.Cells.Item(489, 7).Value = 100103
$ws.Cells.Item(489, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(489, 9).Value = 100103006
$ws.Cells.Item(489, 10).Value = "Nectarín"
$ws.Cells.Item(489, 11).Value = "Early Glo"
$ws.Cells.Item(489, 12).Value = "Primera"
$ws.Cells.Item(489, 13).Value = 200
$ws.Cells.Item(489, 14).Value = 16000
$ws.Cells.Item(489, 15).Value = 16000
$ws.Cells.Item(489, 16).Value = 16000
$ws.Cells.Item(489, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(489, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(489, 19).Value = 889
$ws.Cells.Item(489, 20).Value = 18
